# "Fixed BOM and PnP files"
# In the PnP (Pick-and-Place) sheet, the E9 "Rotation" value was corrected
# from 90 to 0, and the active selection was left on E9 (instead of the
# previous A7:E7 row selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Correct the rotation value for the part in row 9 (column E = "Rotation").
$ws.Range("E9").Value = 0

# Move/leave the selection on the corrected cell, matching the saved file.
$ws.Activate()
$ws.Range("E9").Select()
